# This edit removes the "nafo informado" data row from the municipality
# table. That row (municipio = "nafo informado", casos = 105, obitos = 0)
# is deleted entirely, which shifts every row below it up by one and
# shrinks the used range from A1:C585 to A1:C584.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the row whose column A value is "nafo informado" and delete it
# (shifting the cells below it up), rather than hard-coding row 337, so
# the script is resilient if the sheet layout changes slightly.
$target = $ws.Cells.Find("nafo informado", [System.Reflection.Missing]::Value, -4163, 1)

if ($target -ne $null) {
    $rowNum = $target.Row
} else {
    # Fallback: based on the known layout, the row is 337.
    $rowNum = 337
}

$ws.Rows.Item($rowNum).Delete()
